$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.483.97'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +1.97%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.865.76'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +1.00%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.012'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '311.79'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.13%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4777'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.36%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3748'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +2.27%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07326'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +1.46%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9350'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.87%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +4.93%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07828'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +1.54%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.878.60'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +2.16%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.436'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.550'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +2.16%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '90.31'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +1.71%  '
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -0.14%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008885'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +2.89%  '
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.07%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '27.571.32'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +2.18%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.62'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.50%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.117'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +1.20%  '
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.30%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.75%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '154.47'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +1.27%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +1.46%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.022'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +1.44%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '115.49'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.08898'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.28%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.335'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.33%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +4.02%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.7580'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +1.94%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.613'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +2.71%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.751'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +1.32%  '
$ws.Range('B36').Value = 'TrustWalletToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.119'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +0.15%  '
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02032'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +3.89%  '
$ws.Range('B38').Value = 'Hedera'
$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05262'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -0.10%  '
$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.990'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +0.14%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.5309'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +2.30%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '7.075'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +1.28%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1523'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.92%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.472'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +3.41%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '10.58'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +0.30%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.4801'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +1.56%  '
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +0.11%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '102.88'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +1.35%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.651'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +3.14%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '67.43'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +2.90%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06078'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +0.89%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.9188'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +3.61%  '
